$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 data entries (Week 8 row)
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 0.5
$ws.Range("F12").Value = 3.75
$ws.Range("I12").Value = "Team Project work"

# Column C width was manually resized (no longer auto best-fit)
$ws.Columns.Item(3).ColumnWidth = 10.8

# Active selection moved to K10
$ws.Range("K10").Select()
